$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated team-specific transition-matrix probabilities after adding more simulated games.
$updates = @{
    "B2" = 0.2358078602620087
    "C2" = 0.4759825327510917
    "J2" = 0.02183406113537118
    "P2" = 0.1441048034934498
    "S2" = 0.1222707423580786
    "C3" = 0.01834862385321101
    "J3" = 0.04587155963302753
    "P3" = 0.7064220183486238
    "S3" = 0.2293577981651376
    "J4" = 0.07142857142857142
    "P4" = 0.6071428571428571
    "S4" = 0.3214285714285715
    "B6" = 0.06829268292682927
    "D6" = 0.02439024390243903
    "F6" = 0.02926829268292683
    "J6" = 0.1658536585365854
    "O6" = 0.01463414634146342
    "Q6" = 0.1658536585365854
    "R6" = 0.1170731707317073
    "S6" = 0.4146341463414634
    "B7" = 0.07954545454545454
    "D7" = 0.02840909090909091
    "E7" = 0.005681818181818182
    "F7" = 0.05113636363636364
    "J7" = 0.1704545454545454
    "O7" = 0.01136363636363636
    "Q7" = 0.2159090909090909
    "R7" = 0.0625
    "S7" = 0.375
    "B8" = 0.07223476297968397
    "D8" = 0.01354401805869074
    "F8" = 0.06094808126410835
    "J8" = 0.1038374717832957
    "O8" = 0.01805869074492099
    "Q8" = 0.1693002257336343
    "R8" = 0.0835214446952596
    "S8" = 0.4785553047404063
    "B9" = 0.07851239669421488
    "D9" = 0.02892561983471074
    "F9" = 0.05785123966942149
    "J9" = 0.1115702479338843
    "O9" = 0.01652892561983471
    "Q9" = 0.1859504132231405
    "R9" = 0.1033057851239669
    "S9" = 0.4173553719008264
    "B10" = 0.08451957295373666
    "D10" = 0.02846975088967971
    "E10" = 0.0008896797153024911
    "F10" = 0.07829181494661921
    "J10" = 0.09252669039145907
    "O10" = 0.01779359430604982
    "Q10" = 0.201067615658363
    "R10" = 0.09608540925266904
    "S10" = 0.400355871886121
    "G11" = 0.1259259259259259
    "J11" = 0.0962962962962963
    "K11" = 0.1925925925925926
    "L11" = 0.5703703703703704
    "S11" = 0.01481481481481482
    "G12" = 0.7597402597402597
    "J12" = 0.1883116883116883
    "K12" = 0.01298701298701299
    "L12" = 0.01298701298701299
    "S12" = 0.02597402597402598
    "G13" = 0.75
    "J13" = 0.2222222222222222
    "S13" = 0.02777777777777778
    "F15" = 0.0136986301369863
    "H15" = 0.1552511415525114
    "I15" = 0.1004566210045662
    "J15" = 0.3607305936073059
    "K15" = 0.0593607305936073
    "M15" = 0.0091324200913242
    "O15" = 0.0365296803652968
    "S15" = 0.2648401826484018
    "F16" = 0.0291970802919708
    "H16" = 0.1897810218978102
    "I16" = 0.0948905109489051
    "J16" = 0.3357664233576642
    "K16" = 0.072992700729927
    "M16" = 0.0218978102189781
    "N16" = 0.0072992700729927
    "O16" = 0.08759124087591241
    "S16" = 0.1605839416058394
    "F17" = 0.01913875598086124
    "H17" = 0.1626794258373206
    "I17" = 0.1411483253588517
    "J17" = 0.3803827751196172
    "K17" = 0.08851674641148326
    "M17" = 0.02392344497607655
    "O17" = 0.06220095693779904
    "S17" = 0.1220095693779904
    "F18" = 0.009803921568627451
    "H18" = 0.2009803921568628
    "I18" = 0.09313725490196079
    "J18" = 0.3823529411764706
    "K18" = 0.08823529411764706
    "M18" = 0.02941176470588235
    "O18" = 0.07843137254901961
    "S18" = 0.1176470588235294
    "F19" = 0.01259842519685039
    "H19" = 0.2188976377952756
    "I19" = 0.1007874015748031
    "J19" = 0.3637795275590551
    "K19" = 0.1039370078740157
    "M19" = 0.01181102362204724
    "N19" = 0.001574803149606299
    "O19" = 0.07716535433070866
    "S19" = 0.1094488188976378
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}

